$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC row 6
$wsALC.Range("H6").Value = 10084
$wsALC.Range("I6").Value = 10084
$wsALC.Range("K6").Value = 30252
$wsALC.Range("M6").Value = -30140

# ALC row 58
$wsALC.Range("H58").Value = 1585.1177
$wsALC.Range("J58").Value = 3109.5715
$wsALC.Range("L58").Value = 9328.7145
$wsALC.Range("N58").Value = -9628.7145

# ALC row 74
$wsALC.Range("H74").Value = 3625.25
$wsALC.Range("J74").Value = 3799.8
$wsALC.Range("L74").Value = 3799.8
$wsALC.Range("N74").Value = -5671.8

# ALC row 77
$wsALC.Range("H77").Value = 3625.25
$wsALC.Range("J77").Value = 3799.8
$wsALC.Range("L77").Value = 18999
$wsALC.Range("N77").Value = -28359

# ALC row 100
$wsALC.Range("H100").Value = 13514201
$wsALC.Range("I100").Value = 16129534
$wsALC.Range("K100").Value = 16129534
$wsALC.Range("M100").Value = -16128993

# ALC row 107
$wsALC.Range("H107").Value = 1832.2354
$wsALC.Range("J107").Value = 4805.7144
$wsALC.Range("L107").Value = 4805.7144
$wsALC.Range("N107").Value = -8645.714400000001

# ARM row 32
$wsARM.Range("H32").Value = 7628.1387
$wsARM.Range("I32").Value = 7628.1387
$wsARM.Range("J32").Value = 0
$wsARM.Range("K32").Value = 7628.1387
$wsARM.Range("L32").Value = 0
$wsARM.Range("M32").Value = -7341.1387
$wsARM.Range("N32").ClearContents()

# ARM row 61
$wsARM.Range("H61").Value = 33334512
$wsARM.Range("I61").Value = 38462530
$wsARM.Range("J61").Value = 2403.5
$wsARM.Range("K61").Value = 38462530
$wsARM.Range("L61").Value = 2403.5
$wsARM.Range("M61").Value = -38462318
$wsARM.Range("N61").Value = -2827.5

# ARM row 110
$wsARM.Range("H110").Value = 1052.4
$wsARM.Range("I110").Value = 376.375
$wsARM.Range("K110").Value = 376.375
$wsARM.Range("M110").Value = 1668.625

# ARM row 132
$wsARM.Range("H132").Value = 2731.96
$wsARM.Range("I132").Value = 2041.1562
$wsARM.Range("J132").Value = 3960.0557
$wsARM.Range("K132").Value = 6123.4686
$wsARM.Range("L132").Value = 11880.1671
$wsARM.Range("M132").Value = -3593.4686
$wsARM.Range("N132").Value = -16940.1671

# ARM row 136
$wsARM.Range("H136").Value = 33334512
$wsARM.Range("I136").Value = 38462530
$wsARM.Range("J136").Value = 2403.5
$wsARM.Range("K136").Value = 115387590
$wsARM.Range("L136").Value = 7210.5
$wsARM.Range("M136").Value = -115385040
$wsARM.Range("N136").Value = -12310.5

# BSM row 86
$wsBSM.Range("H86").Value = 2915.577
$wsBSM.Range("I86").Value = 3049.9443
$wsBSM.Range("J86").Value = 2613.25
$wsBSM.Range("K86").Value = 3049.9443
$wsBSM.Range("L86").Value = 2613.25
$wsBSM.Range("M86").Value = -1926.9443
$wsBSM.Range("N86").Value = -4859.25

# BSM row 89
$wsBSM.Range("H89").Value = 2915.577
$wsBSM.Range("I89").Value = 3049.9443
$wsBSM.Range("J89").Value = 2613.25
$wsBSM.Range("K89").Value = 15249.7215
$wsBSM.Range("L89").Value = 13066.25
$wsBSM.Range("M89").Value = -9633.7215
$wsBSM.Range("N89").Value = -24298.25

# BSM row 94
$wsBSM.Range("H94").Value = 13889309
$wsBSM.Range("I94").Value = 14706314
$wsBSM.Range("J94").Value = 225
$wsBSM.Range("K94").Value = 14706314
$wsBSM.Range("L94").Value = 225
$wsBSM.Range("M94").Value = -14705863
$wsBSM.Range("N94").Value = -1127

# CRP row 16
$wsCRP.Range("H16").Value = 66667868
$wsCRP.Range("I16").Value = 90910340
$wsCRP.Range("K16").Value = 90910340
$wsCRP.Range("M16").Value = -90910053

# CRP row 31
$wsCRP.Range("H31").Value = 1265.3771
$wsCRP.Range("I31").Value = 1157.2858
$wsCRP.Range("K31").Value = 1157.2858
$wsCRP.Range("M31").Value = -862.2858000000001

# CRP row 34
$wsCRP.Range("H34").Value = 1265.3771
$wsCRP.Range("I34").Value = 1157.2858
$wsCRP.Range("K34").Value = 1157.2858
$wsCRP.Range("M34").Value = -955.2858000000001

# CRP row 107
$wsCRP.Range("H107").Value = 726.84
$wsCRP.Range("I107").Value = 455.26666
$wsCRP.Range("K107").Value = 455.26666
$wsCRP.Range("M107").Value = 1464.73334

# CRP row 113
$wsCRP.Range("H113").Value = 66667868
$wsCRP.Range("I113").Value = 90910340
$wsCRP.Range("K113").Value = 90910340
$wsCRP.Range("M113").Value = -90908170

# CUL row 6
$wsCUL.Range("H6").Value = 105.333336
$wsCUL.Range("I6").Value = 30.5
$wsCUL.Range("J6").Value = 255
$wsCUL.Range("K6").Value = 91.5
$wsCUL.Range("L6").Value = 765
$wsCUL.Range("M6").Value = 21.5
$wsCUL.Range("N6").Value = -991

# CUL row 11
$wsCUL.Range("H11").Value = 193.57143
$wsCUL.Range("I11").Value = 193.57143
$wsCUL.Range("K11").Value = 580.71429
$wsCUL.Range("M11").Value = -440.71429

# CUL row 12
$wsCUL.Range("H12").Value = 136.3125
$wsCUL.Range("I12").Value = 230
$wsCUL.Range("J12").Value = 93.72727
$wsCUL.Range("K12").Value = 690
$wsCUL.Range("L12").Value = 281.18181
$wsCUL.Range("M12").Value = -517
$wsCUL.Range("N12").Value = -627.18181

# CUL row 61
$wsCUL.Range("H61").Value = 311.5625
$wsCUL.Range("I61").Value = 186
$wsCUL.Range("J61").Value = 520.8333
$wsCUL.Range("K61").Value = 558
$wsCUL.Range("L61").Value = 1562.4999
$wsCUL.Range("M61").Value = -343
$wsCUL.Range("N61").Value = -1992.4999

# CUL row 97
$wsCUL.Range("H97").Value = 1067.1428
$wsCUL.Range("I97").Value = 900
$wsCUL.Range("J97").Value = 1095
$wsCUL.Range("K97").Value = 2700
$wsCUL.Range("L97").Value = 3285
$wsCUL.Range("M97").Value = -2204
$wsCUL.Range("N97").Value = -4277

# CUL row 131
$wsCUL.Range("H131").Value = 24427580
$wsCUL.Range("I131").Value = 83333784
$wsCUL.Range("J131").Value = 52596.516
$wsCUL.Range("K131").Value = 250001352
$wsCUL.Range("L131").Value = 157789.548
$wsCUL.Range("M131").Value = -249996312
$wsCUL.Range("N131").Value = -167869.548

# CUL row 140
$wsCUL.Range("H140").Value = 26768.834
$wsCUL.Range("I140").Value = 60959.94
$wsCUL.Range("K140").Value = 182879.82
$wsCUL.Range("M140").Value = -177699.82

# GSM row 80
$wsGSM.Range("H80").Value = 3260.1765
$wsGSM.Range("I80").Value = 1685
$wsGSM.Range("J80").Value = 4119.364
$wsGSM.Range("K80").Value = 1685
$wsGSM.Range("L80").Value = 4119.364
$wsGSM.Range("M80").Value = -687
$wsGSM.Range("N80").Value = -6115.364

# GSM row 83
$wsGSM.Range("H83").Value = 3260.1765
$wsGSM.Range("I83").Value = 1685
$wsGSM.Range("J83").Value = 4119.364
$wsGSM.Range("K83").Value = 8425
$wsGSM.Range("L83").Value = 20596.82
$wsGSM.Range("M83").Value = -3433
$wsGSM.Range("N83").Value = -30580.82

# GSM row 97
$wsGSM.Range("H97").Value = 984.4
$wsGSM.Range("I97").Value = 950
$wsGSM.Range("J97").Value = 1007.3333
$wsGSM.Range("K97").Value = 950
$wsGSM.Range("L97").Value = 1007.3333
$wsGSM.Range("M97").Value = -454
$wsGSM.Range("N97").Value = -1999.3333

# GSM row 113
$wsGSM.Range("H113").Value = 1161.1875
$wsGSM.Range("J113").Value = 1442.7142
$wsGSM.Range("L113").Value = 1442.7142
$wsGSM.Range("N113").Value = -5782.7142

# GSM row 132
$wsGSM.Range("H132").Value = 3374.3845
$wsGSM.Range("I132").Value = 3048.5264
$wsGSM.Range("J132").Value = 4258.857
$wsGSM.Range("K132").Value = 9145.5792
$wsGSM.Range("L132").Value = 12776.571
$wsGSM.Range("M132").Value = -6615.5792
$wsGSM.Range("N132").Value = -17836.571

# LTW row 40
$wsLTW.Range("H40").Value = 3282
$wsLTW.Range("I40").Value = 2249.5833
$wsLTW.Range("J40").Value = 5759.8
$wsLTW.Range("K40").Value = 2249.5833
$wsLTW.Range("L40").Value = 5759.8
$wsLTW.Range("M40").Value = -2113.5833
$wsLTW.Range("N40").Value = -6031.8

# LTW row 61
$wsLTW.Range("H61").Value = 1192.6364
$wsLTW.Range("I61").Value = 1038.875
$wsLTW.Range("J61").Value = 1602.6666
$wsLTW.Range("K61").Value = 1038.875
$wsLTW.Range("L61").Value = 1602.6666
$wsLTW.Range("M61").Value = -836.875
$wsLTW.Range("N61").Value = -2006.6666

# LTW row 113
$wsLTW.Range("H113").Value = 1192.6364
$wsLTW.Range("I113").Value = 1038.875
$wsLTW.Range("J113").Value = 1602.6666
$wsLTW.Range("K113").Value = 1038.875
$wsLTW.Range("L113").Value = 1602.6666
$wsLTW.Range("M113").Value = 1131.125
$wsLTW.Range("N113").Value = -5942.6666

# WVR row 107
$wsWVR.Range("H107").Value = 585.0968
$wsWVR.Range("I107").Value = 344.64285
$wsWVR.Range("K107").Value = 1033.92855
$wsWVR.Range("M107").Value = 886.0714499999999

# WVR row 113
$wsWVR.Range("H113").Value = 554.2353000000001
$wsWVR.Range("I113").Value = 380.25
$wsWVR.Range("J113").Value = 971.8
$wsWVR.Range("K113").Value = 1140.75
$wsWVR.Range("L113").Value = 2915.4
$wsWVR.Range("M113").Value = 1029.25
$wsWVR.Range("N113").Value = -7255.4
